$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 420.83334
$ws.Range("I32").Value = 338.33334
$ws.Range("K32").Value = 338.33334
$ws.Range("M32").Value = -12.33334000000002
$ws.Range("H113").Value = 3206.923
$ws.Range("I113").Value = 2205.5557
$ws.Range("J113").Value = 3737.0588
$ws.Range("K113").Value = 2205.5557
$ws.Range("L113").Value = 3737.0588
$ws.Range("M113").Value = 1048.4443
$ws.Range("N113").Value = -10245.0588
$ws.Range("H116").Value = 2185.7144
$ws.Range("I116").Value = 1966.6666
$ws.Range("J116").Value = 3500
$ws.Range("K116").Value = 1966.6666
$ws.Range("L116").Value = 3500
$ws.Range("M116").Value = 1475.3334
$ws.Range("N116").Value = -10384
$ws.Range("H132").Value = 4084292.5
$ws.Range("I132").Value = 4764641.5
$ws.Range("J132").Value = 2200
$ws.Range("K132").Value = 14293924.5
$ws.Range("L132").Value = 6600
$ws.Range("M132").Value = -14291394.5
$ws.Range("N132").Value = -11660
$ws.Range("H137").Value = 1201.2241
$ws.Range("I137").Value = 961.4103
$ws.Range("J137").Value = 1693.4736
$ws.Range("K137").Value = 2884.2309
$ws.Range("L137").Value = 5080.4208
$ws.Range("M137").Value = -334.2309
$ws.Range("N137").Value = -10180.4208
$ws.Range("H141").Value = 1809.5
$ws.Range("I141").Value = 678.47455
$ws.Range("J141").Value = 3613.027
$ws.Range("K141").Value = 2035.42365
$ws.Range("L141").Value = 10839.081
$ws.Range("M141").Value = 3144.57635
$ws.Range("N141").Value = -21199.081

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1134.24
$ws.Range("I32").Value = 1139.6364
$ws.Range("J32").Value = 600
$ws.Range("K32").Value = 1139.6364
$ws.Range("L32").Value = 600
$ws.Range("M32").Value = -852.6364000000001
$ws.Range("N32").Value = -1174
$ws.Range("H45").Value = 1405.2
$ws.Range("I45").Value = 1420
$ws.Range("J45").Value = 1400.2667
$ws.Range("K45").Value = 1420
$ws.Range("L45").Value = 1400.2667
$ws.Range("M45").Value = -1043
$ws.Range("N45").Value = -2154.2667
$ws.Range("H61").Value = 1321.125
$ws.Range("I61").Value = 861.26
$ws.Range("J61").Value = 2963.5
$ws.Range("K61").Value = 861.26
$ws.Range("L61").Value = 2963.5
$ws.Range("M61").Value = -649.26
$ws.Range("N61").Value = -3387.5
$ws.Range("H74").Value = 850.093
$ws.Range("I74").Value = 822.3333
$ws.Range("J74").Value = 992.8570999999999
$ws.Range("K74").Value = 822.3333
$ws.Range("L74").Value = 992.8570999999999
$ws.Range("M74").Value = 51.66669999999999
$ws.Range("N74").Value = -2740.8571
$ws.Range("H77").Value = 850.093
$ws.Range("I77").Value = 822.3333
$ws.Range("J77").Value = 992.8570999999999
$ws.Range("K77").Value = 4111.6665
$ws.Range("L77").Value = 4964.2855
$ws.Range("M77").Value = 256.3334999999997
$ws.Range("N77").Value = -13700.2855
$ws.Range("H102").Value = 1207.3334
$ws.Range("I102").Value = 1092.5
$ws.Range("J102").Value = 1666.6666
$ws.Range("K102").Value = 1092.5
$ws.Range("L102").Value = 1666.6666
$ws.Range("M102").Value = 529.5
$ws.Range("N102").Value = -4910.6666
$ws.Range("H118").Value = 40000
$ws.Range("J118").Value = 40000
$ws.Range("L118").Value = 40000
$ws.Range("N118").Value = -43314
$ws.Range("H122").Value = 2122.4
$ws.Range("I122").Value = 1806
$ws.Range("J122").Value = 2333.3333
$ws.Range("K122").Value = 5418
$ws.Range("L122").Value = 6999.999899999999
$ws.Range("M122").Value = -2968
$ws.Range("N122").Value = -11899.9999
$ws.Range("H136").Value = 1321.125
$ws.Range("I136").Value = 861.26
$ws.Range("J136").Value = 2963.5
$ws.Range("K136").Value = 2583.78
$ws.Range("L136").Value = 8890.5
$ws.Range("M136").Value = -33.77999999999975
$ws.Range("N136").Value = -13990.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H116").Value = 28333.334
$ws.Range("J116").Value = 28333.334
$ws.Range("L116").Value = 28333.334
$ws.Range("N116").Value = -37511.334
$ws.Range("H134").Value = 27696.41
$ws.Range("I134").Value = 41894.84
$ws.Range("J134").Value = 2342.0715
$ws.Range("K134").Value = 125684.52
$ws.Range("L134").Value = 7026.2145
$ws.Range("M134").Value = -123149.52
$ws.Range("N134").Value = -12096.2145

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1466.08
$ws.Range("I134").Value = 1433.6364
$ws.Range("J134").Value = 1704
$ws.Range("K134").Value = 4300.9092
$ws.Range("L134").Value = 5112
$ws.Range("M134").Value = -1765.9092
$ws.Range("N134").Value = -10182

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 5636176.5
$ws.Range("I131").Value = 40396
$ws.Range("J131").Value = 7035121.5
$ws.Range("K131").Value = 121188
$ws.Range("L131").Value = 21105364.5
$ws.Range("M131").Value = -116148
$ws.Range("N131").Value = -21115444.5
$ws.Range("H132").Value = 2237.1304
$ws.Range("I132").Value = 1340.3636
$ws.Range("J132").Value = 3059.1667
$ws.Range("K132").Value = 12063.2724
$ws.Range("L132").Value = 27532.5003
$ws.Range("M132").Value = -9533.2724
$ws.Range("N132").Value = -32592.5003

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H120").Value = 37500
$ws.Range("J120").Value = 37500
$ws.Range("L120").Value = 37500
$ws.Range("N120").Value = -47176
$ws.Range("H132").Value = 7357.4473
$ws.Range("I132").Value = 9864.639999999999
$ws.Range("J132").Value = 2535.923
$ws.Range("K132").Value = 29593.92
$ws.Range("L132").Value = 7607.768999999999
$ws.Range("M132").Value = -27063.92
$ws.Range("N132").Value = -12667.769
$ws.Range("H133").Value = 18834.438
$ws.Range("J133").Value = 18834.438
$ws.Range("L133").Value = 18834.438
$ws.Range("N133").Value = -23894.438

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H111").Value = 41500
$ws.Range("J111").Value = 41500
$ws.Range("L111").Value = 41500
$ws.Range("N111").Value = -49680
$ws.Range("H126").Value = 8244.625
$ws.Range("I126").Value = 13883.777
$ws.Range("J126").Value = 994.2857
$ws.Range("K126").Value = 41651.331
$ws.Range("L126").Value = 2982.8571
$ws.Range("M126").Value = -39181.331
$ws.Range("N126").Value = -7922.8571
$ws.Range("H132").Value = 1694.6279
$ws.Range("I132").Value = 1433.3715
$ws.Range("J132").Value = 2837.625
$ws.Range("K132").Value = 4300.1145
$ws.Range("L132").Value = 8512.875
$ws.Range("M132").Value = -1770.1145
$ws.Range("N132").Value = -13572.875
$ws.Range("H136").Value = 3178.1296
$ws.Range("I136").Value = 3550.558
$ws.Range("J136").Value = 1722.2727
$ws.Range("K136").Value = 10651.674
$ws.Range("L136").Value = 5166.8181
$ws.Range("M136").Value = -8101.673999999999
$ws.Range("N136").Value = -10266.8181
